$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Footer "date" field: 30/08/2022 -> 31/08/2022
#    Appears once on the Slide Master and once on every Custom Layout,
#    always on the shape named "Date Placeholder N".
# ---------------------------------------------------------------------
$oldDate = "30/08/2022"
$newDate = "31/08/2022"

$master = $p.SlideMaster

for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $shp = $master.Shapes.Item($si)
    if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 3, "Content Placeholder 2": update the note about the master
#    file, splitting the added sentence into two runs, and reinstate the
#    "master version" sentence as a new run in front of the URL.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$contentShape = $slide3.Shapes.Item("Content Placeholder 2")
$tr = $contentShape.TextFrame.TextRange

$oldNote = "The master version of this file can be found at "
$newNotePart1 = "Note that Microsoft Edge has issues downsizing images, so you may need to manually resize, e.g. using an online tool, for use "
$newNotePart2 = "in websites."

$fullText = $tr.Text
$noteStart0 = $fullText.IndexOf($oldNote)
$noteStart1 = $noteStart0 + 1

# Replace the whole run in one go, then re-select the trailing part so it
# becomes its own run (mirrors the diff: two <a:r> runs).
$tr.Characters($noteStart1, $oldNote.Length).Text = ($newNotePart1 + $newNotePart2)
$tailStart1 = $noteStart1 + $newNotePart1.Length
$tr.Characters($tailStart1, $newNotePart2.Length).Text = $newNotePart2

# Re-add "The master version of this file can be found at " as a new run
# directly in front of the hyperlink text.
$url = "https://github.com/DfE-R-Community/dfe-r-community.github.io/tree/main/resources/logo.pptx"
$fullText2 = $tr.Text
$urlStart0 = $fullText2.IndexOf($url)
$urlStart1 = $urlStart0 + 1

$tr.Characters($urlStart1, $url.Length).Text = ($oldNote + $url)
$tr.Characters($urlStart1, $oldNote.Length).Text = $oldNote
